$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H61").Value = 1838.3334
$ws.Range("I61").Value = 757.5
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 2272.5
$ws.Range("L61").Value = 12000
$ws.Range("M61").Value = -2100.5
$ws.Range("N61").Value = -12344
$ws.Range("H86").Value = 4733.7856
$ws.Range("I86").Value = 2540.3333
$ws.Range("J86").Value = 8682
$ws.Range("K86").Value = 2540.3333
$ws.Range("L86").Value = 8682
$ws.Range("M86").Value = -1417.3333
$ws.Range("N86").Value = -10928
$ws.Range("H89").Value = 4733.7856
$ws.Range("I89").Value = 2540.3333
$ws.Range("J89").Value = 8682
$ws.Range("K89").Value = 12701.6665
$ws.Range("L89").Value = 43410
$ws.Range("M89").Value = -7085.666499999999
$ws.Range("N89").Value = -54642
$ws.Range("H94").Value = 1550
$ws.Range("I94").Value = 1550
$ws.Range("K94").Value = 1550
$ws.Range("M94").Value = -1099
$ws.Range("H106").Value = 102568420
$ws.Range("I106").Value = 33337886
$ws.Range("J106").Value = 333336830
$ws.Range("K106").Value = 33337886
$ws.Range("L106").Value = 333336830
$ws.Range("M106").Value = -33337255
$ws.Range("N106").Value = -333338092
$ws.Range("H107").Value = 16671042
$ws.Range("I107").Value = 20834462
$ws.Range("K107").Value = 20834462
$ws.Range("M107").Value = -20832542

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 80007
$ws.Range("J23").Value = 80007
$ws.Range("L23").Value = 80007
$ws.Range("N23").Value = -80525
$ws.Range("H37").Value = 20000
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
$ws.Range("H102").Value = 7428705
$ws.Range("I102").Value = 7428705
$ws.Range("K102").Value = 7428705
$ws.Range("M102").Value = -7427083
$ws.Range("H105").Value = 38900
$ws.Range("J105").Value = 38900
$ws.Range("L105").Value = 38900
$ws.Range("N105").Value = -45888
$ws.Range("H110").Value = 1150.0588
$ws.Range("I110").Value = 1116.0714
$ws.Range("J110").Value = 1308.6666
$ws.Range("K110").Value = 1116.0714
$ws.Range("L110").Value = 1308.6666
$ws.Range("M110").Value = 928.9286
$ws.Range("N110").Value = -5398.6666
$ws.Range("H122").Value = 1604452.4
$ws.Range("I122").Value = 2332521.8
$ws.Range("J122").Value = 2700
$ws.Range("K122").Value = 6997565.399999999
$ws.Range("L122").Value = 8100
$ws.Range("M122").Value = -6995115.399999999
$ws.Range("N122").Value = -13000

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2794
$ws.Range("I86").Value = 1911.6
$ws.Range("K86").Value = 1911.6
$ws.Range("M86").Value = -788.5999999999999
$ws.Range("H89").Value = 2794
$ws.Range("I89").Value = 1911.6
$ws.Range("K89").Value = 9558
$ws.Range("M89").Value = -3942
$ws.Range("H94").Value = 2106.85
$ws.Range("I94").Value = 1408
$ws.Range("K94").Value = 1408
$ws.Range("M94").Value = -957
$ws.Range("H99").Value = 71430024
$ws.Range("I99").Value = 142858180
$ws.Range("J99").Value = 1870
$ws.Range("K99").Value = 142858180
$ws.Range("L99").Value = 1870
$ws.Range("M99").Value = -142856682
$ws.Range("N99").Value = -4866
$ws.Range("H105").Value = 9877.92
$ws.Range("I105").Value = 14831.267
$ws.Range("J105").Value = 2447.9
$ws.Range("K105").Value = 14831.267
$ws.Range("L105").Value = 2447.9
$ws.Range("M105").Value = -13084.267
$ws.Range("N105").Value = -5941.9
$ws.Range("H107").Value = 1270.2
$ws.Range("I107").Value = 1315.1666
$ws.Range("K107").Value = 1315.1666
$ws.Range("M107").Value = 604.8334
$ws.Range("H138").Value = 54037.145
$ws.Range("J138").Value = 54037.145
$ws.Range("L138").Value = 54037.145
$ws.Range("N138").Value = -64317.145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3089.1667
$ws.Range("I16").Value = 1333.3334
$ws.Range("J16").Value = 4845
$ws.Range("K16").Value = 1333.3334
$ws.Range("L16").Value = 4845
$ws.Range("M16").Value = -1046.3334
$ws.Range("N16").Value = -5419
$ws.Range("H105").Value = 1911.6
$ws.Range("I105").Value = 1181.6666
$ws.Range("J105").Value = 3006.5
$ws.Range("K105").Value = 1181.6666
$ws.Range("L105").Value = 3006.5
$ws.Range("M105").Value = 565.3334
$ws.Range("N105").Value = -6500.5
$ws.Range("H107").Value = 775.41174
$ws.Range("I107").Value = 597.4286
$ws.Range("K107").Value = 597.4286
$ws.Range("M107").Value = 1322.5714
$ws.Range("H113").Value = 3089.1667
$ws.Range("I113").Value = 1333.3334
$ws.Range("J113").Value = 4845
$ws.Range("K113").Value = 1333.3334
$ws.Range("L113").Value = 4845
$ws.Range("M113").Value = 836.6666
$ws.Range("N113").Value = -9185
$ws.Range("H132").Value = 2343.6428
$ws.Range("I132").Value = 2136.7144
$ws.Range("J132").Value = 2964.4285
$ws.Range("K132").Value = 6410.1432
$ws.Range("L132").Value = 8893.2855
$ws.Range("M132").Value = -3880.1432
$ws.Range("N132").Value = -13953.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 6807017
$ws.Range("I75").Value = 3805.2
$ws.Range("J75").Value = 8933021
$ws.Range("K75").Value = 11415.6
$ws.Range("L75").Value = 26799063
$ws.Range("M75").Value = -10417.6
$ws.Range("N75").Value = -26801059
$ws.Range("H78").Value = 6807017
$ws.Range("I78").Value = 3805.2
$ws.Range("J78").Value = 8933021
$ws.Range("K78").Value = 34246.8
$ws.Range("L78").Value = 80397189
$ws.Range("M78").Value = -29254.8
$ws.Range("N78").Value = -80407173

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H82").Value = 34000
$ws.Range("J82").Value = 34000
$ws.Range("L82").Value = 34000
$ws.Range("N82").Value = -34766
$ws.Range("H85").Value = 34000
$ws.Range("J85").Value = 34000
$ws.Range("L85").Value = 34000
$ws.Range("N85").Value = -36652
$ws.Range("H113").Value = 76924130
$ws.Range("I113").Value = 100000950
$ws.Range("J113").Value = 1366.6666
$ws.Range("K113").Value = 100000950
$ws.Range("L113").Value = 1366.6666
$ws.Range("M113").Value = -99998780
$ws.Range("N113").Value = -5706.6666
$ws.Range("H122").Value = 70989280
$ws.Range("I122").Value = 76059730
$ws.Range("K122").Value = 228179190
$ws.Range("M122").Value = -228176740
$ws.Range("H132").Value = 3405.5
$ws.Range("I132").Value = 3114.6428
$ws.Range("K132").Value = 9343.928400000001
$ws.Range("M132").Value = -6813.928400000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2954.6667
$ws.Range("I61").Value = 3074.25
$ws.Range("J61").Value = 1998
$ws.Range("K61").Value = 3074.25
$ws.Range("L61").Value = 1998
$ws.Range("M61").Value = -2872.25
$ws.Range("N61").Value = -2402
$ws.Range("H82").Value = 13205.385
$ws.Range("I82").Value = 17850.166
$ws.Range("J82").Value = 9224.143
$ws.Range("K82").Value = 17850.166
$ws.Range("L82").Value = 9224.143
$ws.Range("M82").Value = -17489.166
$ws.Range("N82").Value = -9946.143
$ws.Range("H85").Value = 13205.385
$ws.Range("I85").Value = 17850.166
$ws.Range("J85").Value = 9224.143
$ws.Range("K85").Value = 17850.166
$ws.Range("L85").Value = 9224.143
$ws.Range("M85").Value = -16602.166
$ws.Range("N85").Value = -11720.143
$ws.Range("H93").Value = 100042856
$ws.Range("I93").Value = 200000
$ws.Range("J93").Value = 125003576
$ws.Range("K93").Value = 200000
$ws.Range("L93").Value = 125003576
$ws.Range("M93").Value = -198752
$ws.Range("N93").Value = -125006072
$ws.Range("H100").Value = 2722.25
$ws.Range("I100").Value = 2380
$ws.Range("J100").Value = 3292.6667
$ws.Range("K100").Value = 2380
$ws.Range("L100").Value = 3292.6667
$ws.Range("M100").Value = -1839
$ws.Range("N100").Value = -4374.6667
$ws.Range("H113").Value = 2954.6667
$ws.Range("I113").Value = 3074.25
$ws.Range("J113").Value = 1998
$ws.Range("K113").Value = 3074.25
$ws.Range("L113").Value = 1998
$ws.Range("M113").Value = -904.25
$ws.Range("N113").Value = -6338
$ws.Range("H122").Value = 6787921
$ws.Range("I122").Value = 7145305
$ws.Range("J122").Value = 5001000
$ws.Range("K122").Value = 21435915
$ws.Range("L122").Value = 15003000
$ws.Range("M122").Value = -21433465
$ws.Range("N122").Value = -15007900
$ws.Range("H133").Value = 25913
$ws.Range("J133").Value = 25913
$ws.Range("L133").Value = 25913
$ws.Range("N133").Value = -30973
$ws.Range("H136").Value = 3703.3262
$ws.Range("I136").Value = 2733.2122
$ws.Range("J136").Value = 6165.923
$ws.Range("K136").Value = 8199.6366
$ws.Range("L136").Value = 18497.769
$ws.Range("M136").Value = -5649.6366
$ws.Range("N136").Value = -23597.769

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 76923710
$ws.Range("I107").Value = 142857680
$ws.Range("J107").Value = 742.1667
$ws.Range("K107").Value = 428573040
$ws.Range("L107").Value = 2226.5001
$ws.Range("M107").Value = -428571120
$ws.Range("N107").Value = -6066.5001
$ws.Range("H132").Value = 1536.7213
$ws.Range("I132").Value = 1238.1777
$ws.Range("J132").Value = 2376.375
$ws.Range("K132").Value = 3714.5331
$ws.Range("L132").Value = 7129.125
$ws.Range("M132").Value = -1184.5331
$ws.Range("N132").Value = -12189.125
$ws.Range("H136").Value = 1265.0244
$ws.Range("I136").Value = 761.8889
$ws.Range("J136").Value = 2235.3572
$ws.Range("K136").Value = 2285.6667
$ws.Range("L136").Value = 6706.071599999999
$ws.Range("M136").Value = 264.3332999999998
$ws.Range("N136").Value = -11806.0716
